# Atualizado por script em 21-12-2023 20:46
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, [int]$row1, [int]$row2)
    # Columns F (6) through V (22) hold the match data; swap them between the two rows.
    for ($col = 6; $col -le 22; $col++) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value2
        $v2 = $c2.Value2
        $c1.Value2 = $v2
        $c2.Value2 = $v1
    }
}

# Reorder the matches that were shuffled in the source feed.
Swap-RowData $ws 124 126
Swap-RowData $ws 125 127
Swap-RowData $ws 145 147
Swap-RowData $ws 163 164
Swap-RowData $ws 168 169

# Append the newly scraped match as row 171, copying formatting from row 170.
$ws.Range("A170:V170").Copy($ws.Range("A171:V171"))

$ws.Cells.Item(171, 1).Value2 = 170
$ws.Cells.Item(171, 2).Value2 = "england"
$ws.Cells.Item(171, 3).Value2 = "premier-league"
$ws.Cells.Item(171, 4).Value2 = "2023-2024"
$ws.Cells.Item(171, 5).Value2 = 45281.875
$ws.Cells.Item(171, 6).Value2 = "Crystal Palace"
$ws.Cells.Item(171, 7).Value2 = 1
$ws.Cells.Item(171, 8).Value2 = "Brighton"
$ws.Cells.Item(171, 9).Value2 = 1
$ws.Cells.Item(171, 10).Value2 = 2.98
$ws.Cells.Item(171, 11).Value2 = "09/12/2023 00:02"
$ws.Cells.Item(171, 12).Value2 = 3.23
$ws.Cells.Item(171, 13).Value2 = "21/12/2023 20:53"
$ws.Cells.Item(171, 14).Value2 = 3.42
$ws.Cells.Item(171, 15).Value2 = "09/12/2023 00:02"
$ws.Cells.Item(171, 16).Value2 = 3.37
$ws.Cells.Item(171, 17).Value2 = "21/12/2023 20:47"
$ws.Cells.Item(171, 18).Value2 = 2.3
$ws.Cells.Item(171, 19).Value2 = "09/12/2023 00:02"
$ws.Cells.Item(171, 20).Value2 = 2.38
$ws.Cells.Item(171, 21).Value2 = "21/12/2023 20:53"
$ws.Cells.Item(171, 22).Value2 = "https://www.betexplorer.com/football/england/premier-league/crystal-palace-brighton/4AtrLXsT/"

Write-Host "Done"
